$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (shifts old Z.. right by one, making room
# for the new "Remarks Visual Inspection" column).
$ws.Range("Z1").EntireColumn.Insert()

# Rename the "Remarks" header (still in column Y) to be specific to the load test,
# and move its old value (the annual visual-inspection remark) into the new
# column Z, giving it its own "Remarks Visual Inspection" header.
$ws.Range("Y1").Value = "Remarks Load Test"
$ws.Range("Z1").Value = "Remarks Visual Inspection"

$ws.Range("Z2").Value = $ws.Range("Y2").Value2
$ws.Range("Y2").Value = "Load test Carried out by Master & DNV Inspector in Durban 26.11.2022"

# Update the view state (scrolled/selected cell) to match the rework.
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("Y3").Select()

Write-Host "done"
